# Update "Pais" worksheet: refresh COVID-19 country stats and fix
# country row ordering (alphabetical-ish swaps caused by new totals),
# plus bump the "Datos actualizados" timestamp footer.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Footer timestamp (row 1)
$ws.Cells.Item(1,1).Value2 = "Datos actualizados a 24 de Octubre de 2020 a las 17:36"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Cells.Item(4,2).Value2 = 8764777
$ws.Cells.Item(4,3).Value2 = 17824
$ws.Cells.Item(4,4).Value2 = 5703372
$ws.Cells.Item(4,5).Value2 = 2831908
$ws.Cells.Item(4,7).Value2 = 213
$ws.Cells.Item(4,8).Value2 = 229497

# Row 5: India -> India
$ws.Cells.Item(5,2).Value2 = 7829226
$ws.Cells.Item(5,3).Value2 = 15558
$ws.Cells.Item(5,4).Value2 = 7030903
$ws.Cells.Item(5,5).Value2 = 680087
$ws.Cells.Item(5,7).Value2 = 244
$ws.Cells.Item(5,8).Value2 = 118236

# Row 17: Chile -> Italia
$ws.Cells.Item(17,1).Value2 = "Italia"
$ws.Cells.Item(17,2).Value2 = 504509
$ws.Cells.Item(17,3).Value2 = 19644
$ws.Cells.Item(17,4).Value2 = 264117
$ws.Cells.Item(17,5).Value2 = 203182
$ws.Cells.Item(17,7).Value2 = 151
$ws.Cells.Item(17,8).Value2 = 37210

# Row 18: Italia -> Chile
$ws.Cells.Item(18,1).Value2 = "Chile"
$ws.Cells.Item(18,2).Value2 = 500542
$ws.Cells.Item(18,3).Value2 = 1631
$ws.Cells.Item(18,4).Value2 = 476576
$ws.Cells.Item(18,5).Value2 = 10074
$ws.Cells.Item(18,7).Value2 = 48
$ws.Cells.Item(18,8).Value2 = 13892

# Row 20: Alemania -> Alemania
$ws.Cells.Item(20,2).Value2 = 422506
$ws.Cells.Item(20,3).Value2 = 5156
$ws.Cells.Item(20,5).Value2 = 98311
$ws.Cells.Item(20,7).Value2 = 5
$ws.Cells.Item(20,8).Value2 = 10095

# Row 30: Paises Bajos -> Paises Bajos
$ws.Cells.Item(30,2).Value2 = 281052
$ws.Cells.Item(30,3).Value2 = 8651
$ws.Cells.Item(30,7).Value2 = 55
$ws.Cells.Item(30,8).Value2 = 7019

# Row 33: Canada -> Canada
$ws.Cells.Item(33,2).Value2 = 212713
$ws.Cells.Item(33,3).Value2 = 981
$ws.Cells.Item(33,4).Value2 = 178504
$ws.Cells.Item(33,5).Value2 = 24315
$ws.Cells.Item(33,7).Value2 = 6
$ws.Cells.Item(33,8).Value2 = 9894

# Row 41: Emiratos Arabes Unidos -> Republica Dominicana
$ws.Cells.Item(41,1).Value2 = "Republica Dominicana"
$ws.Cells.Item(41,2).Value2 = 124018
$ws.Cells.Item(41,3).Value2 = 634
$ws.Cells.Item(41,4).Value2 = 102353
$ws.Cells.Item(41,5).Value2 = 19445
$ws.Cells.Item(41,7).Value2 = 6
$ws.Cells.Item(41,8).Value2 = 2220

# Row 42: Republica Dominicana -> Emiratos Arabes Unidos
$ws.Cells.Item(42,1).Value2 = "Emiratos Arabes Unidos"
$ws.Cells.Item(42,2).Value2 = 123764
$ws.Cells.Item(42,3).Value2 = 1491
$ws.Cells.Item(42,4).Value2 = 116894
$ws.Cells.Item(42,5).Value2 = 6395
$ws.Cells.Item(42,8).Value2 = 475

# Row 61: Moldavia -> Moldavia
$ws.Cells.Item(61,2).Value2 = 71089
$ws.Cells.Item(61,3).Value2 = 833
$ws.Cells.Item(61,4).Value2 = 51719
$ws.Cells.Item(61,5).Value2 = 17701
$ws.Cells.Item(61,7).Value2 = 15
$ws.Cells.Item(61,8).Value2 = 1669

# Row 63: Uzbekistan -> Uzbekistan
$ws.Cells.Item(63,2).Value2 = 64923
$ws.Cells.Item(63,3).Value2 = 199
$ws.Cells.Item(63,5).Value2 = 2346
$ws.Cells.Item(63,7).Value2 = 2
$ws.Cells.Item(63,8).Value2 = 544

# Row 66: Singapur -> Singapur
$ws.Cells.Item(66,4).Value2 = 57844
$ws.Cells.Item(66,5).Value2 = 93

# Row 72: Estado de Palestina -> Jordania
$ws.Cells.Item(72,1).Value2 = "Jordania"
$ws.Cells.Item(72,2).Value2 = 50750
$ws.Cells.Item(72,3).Value2 = 1820
$ws.Cells.Item(72,4).Value2 = 7508
$ws.Cells.Item(72,5).Value2 = 42702
$ws.Cells.Item(72,7).Value2 = 32
$ws.Cells.Item(72,8).Value2 = 540

# Row 73: Azerbaiyan -> Estado de Palestina
$ws.Cells.Item(73,1).Value2 = "Estado de Palestina"
$ws.Cells.Item(73,2).Value2 = 49989
$ws.Cells.Item(73,3).Value2 = 410
$ws.Cells.Item(73,4).Value2 = 43232
$ws.Cells.Item(73,5).Value2 = 6314
$ws.Cells.Item(73,7).Value2 = 4
$ws.Cells.Item(73,8).Value2 = 443

# Row 74: Jordania -> Azerbaiyan
$ws.Cells.Item(74,1).Value2 = "Azerbaiyan"
$ws.Cells.Item(74,2).Value2 = 49013
$ws.Cells.Item(74,3).Value2 = 792
$ws.Cells.Item(74,4).Value2 = 41051
$ws.Cells.Item(74,5).Value2 = 7298
$ws.Cells.Item(74,7).Value2 = 8
$ws.Cells.Item(74,8).Value2 = 664

# Row 81: Dinamarca -> Bosnia y Herzegovina
$ws.Cells.Item(81,1).Value2 = "Bosnia y Herzegovina"
$ws.Cells.Item(81,2).Value2 = 39758
$ws.Cells.Item(81,3).Value2 = 1265
$ws.Cells.Item(81,4).Value2 = 26368
$ws.Cells.Item(81,5).Value2 = 12315
$ws.Cells.Item(81,7).Value2 = 10
$ws.Cells.Item(81,8).Value2 = 1075

# Row 82: Serbia -> Dinamarca
$ws.Cells.Item(82,1).Value2 = "Dinamarca"
$ws.Cells.Item(82,2).Value2 = 39411
$ws.Cells.Item(82,3).Value2 = 789
$ws.Cells.Item(82,4).Value2 = 31701
$ws.Cells.Item(82,5).Value2 = 7010
$ws.Cells.Item(82,8).Value2 = 700

# Row 83: Bosnia y Herzegovina -> Serbia
$ws.Cells.Item(83,1).Value2 = "Serbia"
$ws.Cells.Item(83,2).Value2 = 38872
$ws.Cells.Item(83,3).Value2 = 757
$ws.Cells.Item(83,4).Value2 = 31536
$ws.Cells.Item(83,5).Value2 = 6547
$ws.Cells.Item(83,7).Value2 = 3
$ws.Cells.Item(83,8).Value2 = 789

# Row 96: Albania -> Albania
$ws.Cells.Item(96,2).Value2 = 18858
$ws.Cells.Item(96,3).Value2 = 302
$ws.Cells.Item(96,4).Value2 = 10548
$ws.Cells.Item(96,5).Value2 = 7837
$ws.Cells.Item(96,7).Value2 = 4
$ws.Cells.Item(96,8).Value2 = 473

# Row 99: Montenegro -> Montenegro
$ws.Cells.Item(99,2).Value2 = 16629
$ws.Cells.Item(99,3).Value2 = 193
$ws.Cells.Item(99,4).Value2 = 12601
$ws.Cells.Item(99,5).Value2 = 3765
$ws.Cells.Item(99,7).Value2 = 8
$ws.Cells.Item(99,8).Value2 = 263

# Row 104: Luxemburgo -> Luxemburgo
$ws.Cells.Item(104,2).Value2 = 13713
$ws.Cells.Item(104,3).Value2 = 862
$ws.Cells.Item(104,5).Value2 = 4484
$ws.Cells.Item(104,7).Value2 = 3
$ws.Cells.Item(104,8).Value2 = 144

# Row 117: Jamaica -> Jamaica
$ws.Cells.Item(117,2).Value2 = 8670
$ws.Cells.Item(117,3).Value2 = 32
$ws.Cells.Item(117,4).Value2 = 4209
$ws.Cells.Item(117,5).Value2 = 4275
$ws.Cells.Item(117,7).Value2 = 4
$ws.Cells.Item(117,8).Value2 = 186

# Row 123: Cuba -> Cuba
$ws.Cells.Item(123,2).Value2 = 6534
$ws.Cells.Item(123,3).Value2 = 55
$ws.Cells.Item(123,4).Value2 = 5927
$ws.Cells.Item(123,5).Value2 = 479

# Row 177: Burundi -> Burundi
$ws.Cells.Item(177,2).Value2 = 555
$ws.Cells.Item(177,3).Value2 = 2
$ws.Cells.Item(177,5).Value2 = 57

# Row 186: Mongolia -> Liechtenstein
$ws.Cells.Item(186,1).Value2 = "Liechtenstein"
$ws.Cells.Item(186,2).Value2 = 340
$ws.Cells.Item(186,3).Value2 = 16
$ws.Cells.Item(186,4).Value2 = 170
$ws.Cells.Item(186,5).Value2 = 169
$ws.Cells.Item(186,8).Value2 = 1

# Row 187: Butan -> Mongolia
$ws.Cells.Item(187,1).Value2 = "Mongolia"
$ws.Cells.Item(187,2).Value2 = 337
$ws.Cells.Item(187,3).Value2 = 9
$ws.Cells.Item(187,4).Value2 = 312
$ws.Cells.Item(187,5).Value2 = 25

# Row 188: Liechtenstein -> Butan
$ws.Cells.Item(188,1).Value2 = "Butan"
$ws.Cells.Item(188,2).Value2 = 336
$ws.Cells.Item(188,4).Value2 = 306
$ws.Cells.Item(188,5).Value2 = 30
$ws.Cells.Item(188,8).Value2 = 0

# Row 216: Montserrat -> Islas Malvinas
$ws.Cells.Item(216,1).Value2 = "Islas Malvinas"
$ws.Cells.Item(216,4).Value2 = 13
$ws.Cells.Item(216,8).Value2 = 0

# Row 217: Islas Malvinas -> Montserrat
$ws.Cells.Item(217,1).Value2 = "Montserrat"
$ws.Cells.Item(217,4).Value2 = 12
$ws.Cells.Item(217,8).Value2 = 1

